$d = $word.ActiveDocument

# --- 1. First paragraph: "This is a Microsoft word document." ---
# Append two trailing spaces to the existing (uncolored) run.
$d.Content.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2)

$p1 = $d.Paragraphs(1)

# Insert the 3 new colored runs right before the paragraph mark, one at a
# time, so each InsertAfter call produces its own <w:r> (matching the
# target which keeps 3 separate runs instead of merging them).
$insertPos = $p1.Range.End - 1
$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter([char]0x0028 + "This is a change " + [char]0x2013 + " Ve")
$r1.Font.Color = 192

$insertPos = $p1.Range.End - 1
$r2 = $d.Range($insertPos, $insertPos)
$r2.InsertAfter("rsion for branch alternate")
$r2.Font.Color = 192

$insertPos = $p1.Range.End - 1
$r3 = $d.Range($insertPos, $insertPos)
$r3.InsertAfter(")")
$r3.Font.Color = 192

# --- 2. New shaded empty paragraph after the "Free at last..." paragraph ---
# Using raw OOXML insertion gives an exact, un-inherited <w:p><w:pPr><w:shd .../></w:pPr></w:p>
# (the Paragraphs/InsertParagraphAfter object-model path always materializes
# the predecessor's inherited run/paragraph-mark formatting as an explicit run).
$endRange = $d.Range($d.Content.End, $d.Content.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)
